$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 13 new rows above row 2, shifting existing data down
$ws.Range("A2:C14").EntireRow.Insert()

# New data for the 13 inserted rows
$newData = @(
    @(-0.02587450614997311, 0.008312130346894184, -0.03830997752291809),
    @(-0.03166318188110984, 0.005548692618807229, 0.004886921495199072),
    @(-0.01120646846746765, -0.01967131506119446, 0.04684027514996975),
    @(-0.005999569415247833, -0.01865320652723307, 0.03610649971025327),
    @(-0.01938042674391034, 0.02639810448246341, -0.02997603196473341),
    @(0.04876013357369639, 0.03156137033476548, -0.02466732121649261),
    @(0.07661995095466947, 0.006654067996091003, 0.02539454134447233),
    @(-0.008406669134274104, 0.01747510954737656, 0.02251474718962388),
    @(-0.01705332100391391, 0.02540181328852965, 0.02260201397196701),
    @(0.004734205614243124, 0.03287764106478005, -0.002007128720703896),
    @(0.002516182849094959, 0.03615013385812434, -0.01502437632353524),
    @(-0.01514073128678966, 0.0402807449655873, -0.0335248665263255),
    @(0.004646939025925644, 0.05990842623370023, -0.04932736818279529)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}

# Remove formatting inherited from the insert so new cells have the default style
$ws.Range("A2:C14").ClearFormats()

# Delete the trailing 3 rows that fell outside the new range (old rows 19-21, now at 32-34)
$ws.Range("A32:C34").EntireRow.Delete()
